$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for the "Neutrophils" sending cluster (rows 14-17 in the old 1-based layout),
# whose target-cluster breakdown (ECs/FAPs/MuSCs/Resolving-Mac) is no longer part of the refreshed TPM output.
$ws.Rows("14:17").Delete()

# Refresh all remaining data rows (2-13) with the newly computed TPM-derived NATMI metrics.
# "Resolving-Mac" is no longer used as a Target cluster value anywhere, so it naturally drops out
# of the shared-strings table once every reference to it below is replaced.
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 21.23829066666667
$ws.Range("H2").Value = 63.714872
$ws.Range("I2").Value = 0.9042366413687101
$ws.Range("J2").Value = 0.90423664136871
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.114581
$ws.Range("N2").Value = 0.343743
$ws.Range("O2").Value = 0.007635610029470834
$ws.Range("P2").Value = 0.007635610029470834
$ws.Range("Q2").Value = 2.433504582877334
$ws.Range("R2").Value = 21.901541245896
$ws.Range("S2").Value = 0.006904398367849944
$ws.Range("T2").Value = 0.006904398367849943

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 21.23829066666667
$ws.Range("H3").Value = 63.714872
$ws.Range("I3").Value = 0.9042366413687101
$ws.Range("J3").Value = 0.90423664136871
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.85444966666667
$ws.Range("N3").Value = 44.563349
$ws.Range("O3").Value = 0.9898917347297518
$ws.Range("P3").Value = 0.9898917347297518
$ws.Range("Q3").Value = 315.4831197140365
$ws.Range("R3").Value = 2839.348077426328
$ws.Range("S3").Value = 0.8950963775306768
$ws.Range("T3").Value = 0.8950963775306767

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 21.23829066666667
$ws.Range("H4").Value = 63.714872
$ws.Range("I4").Value = 0.9042366413687101
$ws.Range("J4").Value = 0.90423664136871
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.037105
$ws.Range("N4").Value = 0.111315
$ws.Range("O4").Value = 0.0024726552407774
$ws.Range("P4").Value = 0.0024726552407774
$ws.Range("Q4").Value = 0.7880467751866667
$ws.Range("R4").Value = 7.09242097668
$ws.Range("S4").Value = 0.002235865470183295
$ws.Range("T4").Value = 0.002235865470183295

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.394651
$ws.Range("H5").Value = 4.183953000000001
$ws.Range("I5").Value = 0.05937834432696559
$ws.Range("J5").Value = 0.05937834432696559
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.114581
$ws.Range("N5").Value = 0.343743
$ws.Range("O5").Value = 0.007635610029470834
$ws.Range("P5").Value = 0.007635610029470834
$ws.Range("Q5").Value = 0.159800506231
$ws.Range("R5").Value = 1.438204556079
$ws.Range("S5").Value = 0.0004533898814763511
$ws.Range("T5").Value = 0.0004533898814763511

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.394651
$ws.Range("H6").Value = 4.183953000000001
$ws.Range("I6").Value = 0.05937834432696559
$ws.Range("J6").Value = 0.05937834432696559
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.85444966666667
$ws.Range("N6").Value = 44.563349
$ws.Range("O6").Value = 0.9898917347297518
$ws.Range("P6").Value = 0.9898917347297518
$ws.Range("Q6").Value = 20.71677308206634
$ws.Range("R6").Value = 186.4509577385971
$ws.Range("S6").Value = 0.05877813227120049
$ws.Range("T6").Value = 0.05877813227120049

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.394651
$ws.Range("H7").Value = 4.183953000000001
$ws.Range("I7").Value = 0.05937834432696559
$ws.Range("J7").Value = 0.05937834432696559
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.037105
$ws.Range("N7").Value = 0.111315
$ws.Range("O7").Value = 0.0024726552407774
$ws.Range("P7").Value = 0.0024726552407774
$ws.Range("Q7").Value = 0.05174852535500001
$ws.Range("R7").Value = 0.4657367281950001
$ws.Range("S7").Value = 0.0001468221742887565
$ws.Range("T7").Value = 0.0001468221742887565

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6246503333333333
$ws.Range("H8").Value = 1.873951
$ws.Range("I8").Value = 0.02659497076804196
$ws.Range("J8").Value = 0.02659497076804196
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.114581
$ws.Range("N8").Value = 0.343743
$ws.Range("O8").Value = 0.007635610029470834
$ws.Range("P8").Value = 0.007635610029470834
$ws.Range("Q8").Value = 0.07157305984366666
$ws.Range("R8").Value = 0.6441575385930001
$ws.Range("S8").Value = 0.0002030688255299448
$ws.Range("T8").Value = 0.0002030688255299448

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6246503333333333
$ws.Range("H9").Value = 1.873951
$ws.Range("I9").Value = 0.02659497076804196
$ws.Range("J9").Value = 0.02659497076804196
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 14.85444966666667
$ws.Range("N9").Value = 44.563349
$ws.Range("O9").Value = 0.9898917347297518
$ws.Range("P9").Value = 0.9898917347297518
$ws.Range("Q9").Value = 9.278836935766556
$ws.Range("R9").Value = 83.509532421899
$ws.Range("S9").Value = 0.0263261417486641
$ws.Range("T9").Value = 0.0263261417486641

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6246503333333333
$ws.Range("H10").Value = 1.873951
$ws.Range("I10").Value = 0.02659497076804196
$ws.Range("J10").Value = 0.02659497076804196
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.037105
$ws.Range("N10").Value = 0.111315
$ws.Range("O10").Value = 0.0024726552407774
$ws.Range("P10").Value = 0.0024726552407774
$ws.Range("Q10").Value = 0.02317765061833333
$ws.Range("R10").Value = 0.208598855565
$ws.Range("S10").Value = [double]"6.576019384792071E-05"
$ws.Range("T10").Value = [double]"6.576019384792071E-05"

$ws.Range("A11").Value = "Neutrophils"
$ws.Range("B11").Value = "Efna1"
$ws.Range("C11").Value = "Epha3"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.229944
$ws.Range("H11").Value = 0.689832
$ws.Range("I11").Value = 0.009790043536282392
$ws.Range("J11").Value = 0.009790043536282392
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.114581
$ws.Range("N11").Value = 0.343743
$ws.Range("O11").Value = 0.007635610029470834
$ws.Range("P11").Value = 0.007635610029470834
$ws.Range("Q11").Value = 0.026347213464
$ws.Range("R11").Value = 0.237124921176
$ws.Range("S11").Value = [double]"7.475295461459394E-05"
$ws.Range("T11").Value = [double]"7.475295461459394E-05"

$ws.Range("A12").Value = "Neutrophils"
$ws.Range("B12").Value = "Efna1"
$ws.Range("C12").Value = "Epha3"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.229944
$ws.Range("H12").Value = 0.689832
$ws.Range("I12").Value = 0.009790043536282392
$ws.Range("J12").Value = 0.009790043536282392
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 14.85444966666667
$ws.Range("N12").Value = 44.563349
$ws.Range("O12").Value = 0.9898917347297518
$ws.Range("P12").Value = 0.9898917347297518
$ws.Range("Q12").Value = 3.415691574152
$ws.Range("R12").Value = 30.741224167368
$ws.Range("S12").Value = 0.00969108317921037
$ws.Range("T12").Value = 0.00969108317921037

$ws.Range("A13").Value = "Neutrophils"
$ws.Range("B13").Value = "Efna1"
$ws.Range("C13").Value = "Epha3"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.229944
$ws.Range("H13").Value = 0.689832
$ws.Range("I13").Value = 0.009790043536282392
$ws.Range("J13").Value = 0.009790043536282392
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.037105
$ws.Range("N13").Value = 0.111315
$ws.Range("O13").Value = 0.0024726552407774
$ws.Range("P13").Value = 0.0024726552407774
$ws.Range("Q13").Value = 0.008532072119999999
$ws.Range("R13").Value = 0.07678864907999999
$ws.Range("S13").Value = [double]"2.420740245742757E-05"
$ws.Range("T13").Value = [double]"2.420740245742757E-05"
